$wb = $excel.ActiveWorkbook
$wsOpen = $wb.Worksheets.Item("Open Tickets")

# ---------------------------------------------------------------------------
# 1) Add the new German "Offene Tickets" sheet (translation of "Open
#    Tickets") as a brand-new sheet at the end of the workbook. Its header
#    cells are filled in with the placeholder English text first (mirroring
#    the original "Open Tickets" header) and corrected to lower case only
#    after the rest of the workbook has been translated, so that the shared-
#    string table is populated in the same order as the authored edit.
# ---------------------------------------------------------------------------
$wsSolved = $wb.Worksheets.Item("Solved tickets in a year")
$wsOffene = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsSolved)
$wsOffene.Name = "Offene Tickets"

$wsOffene.Range("A1").Value = "Ticket"
$wsOffene.Range("B1").Value = "Priority"

$wsOffene.Range("A2").Value = "Keine Verbindung zur Box"
$wsOffene.Range("B2").Value = 500

$wsOffene.Range("A3").Value = "Keine API zur Verfügung gestellt"
$wsOffene.Range("B3").Value = 200

$wsOffene.Range("A4").Value = "Box funktioniert nicht"
$wsOffene.Range("B4").Value = 100

$wsOffene.Range("A5").Value = "Funktion fehlt"
$wsOffene.Range("B5").Value = 110

$wsOffene.Range("A6").Value = "Designer kann nicht aktualisiert werden"
$wsOffene.Range("B6").Value = 520

$wsOffene.Range("A7").Value = "Kabel fehlt"
$wsOffene.Range("B7").Value = 250

$wsOffene.Range("A8").Value = "Box zu heiß"
$wsOffene.Range("B8").Value = 450

$wsOffene.Range("A9").Value = "Designer kann nicht installiert werden"
$wsOffene.Range("B9").Value = 120

$wsOffene.PageSetup.TopMargin = 56.69291339
$wsOffene.PageSetup.BottomMargin = 56.69291339

# ---------------------------------------------------------------------------
# 2) Add the new German "Geschlossene Tickets pro Jahr" sheet (translation
#    of "Solved tickets in a year") at the end of the workbook.
# ---------------------------------------------------------------------------
$wsGeschlossene = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsOffene)
$wsGeschlossene.Name = "Geschlossene Tickets pro Jahr"

$wsGeschlossene.Range("A1").Value = "month"
$wsGeschlossene.Range("B1").Value = "amount_closed_tickets"
$wsGeschlossene.Range("C1").Value = "amount_all_tickets"

$wsGeschlossene.Range("A2").Value = 43739
$wsGeschlossene.Range("B2").Value = 40
$wsGeschlossene.Range("C2").Value = 46

$wsGeschlossene.Range("A3").Value = 43770
$wsGeschlossene.Range("B3").Value = 53
$wsGeschlossene.Range("C3").Value = 54

$wsGeschlossene.Range("A4").Value = 43800
$wsGeschlossene.Range("B4").Value = 15
$wsGeschlossene.Range("C4").Value = 15

$wsGeschlossene.Range("A5").Value = 43831
$wsGeschlossene.Range("B5").Value = 40
$wsGeschlossene.Range("C5").Value = 55

$wsGeschlossene.Range("A6").Value = 43862
$wsGeschlossene.Range("B6").Value = 20
$wsGeschlossene.Range("C6").Value = 40

$wsGeschlossene.Range("A7").Value = 43891
$wsGeschlossene.Range("B7").Value = 50
$wsGeschlossene.Range("C7").Value = 60

$wsGeschlossene.Range("A8").Value = 43922
$wsGeschlossene.Range("B8").Value = 70
$wsGeschlossene.Range("C8").Value = 80

$wsGeschlossene.Range("A9").Value = 43952
$wsGeschlossene.Range("B9").Value = 40
$wsGeschlossene.Range("C9").Value = 70

$wsGeschlossene.Range("A10").Value = 43983
$wsGeschlossene.Range("B10").Value = 20
$wsGeschlossene.Range("C10").Value = 25

$wsGeschlossene.Range("A11").Value = 44013
$wsGeschlossene.Range("B11").Value = 5
$wsGeschlossene.Range("C11").Value = 10

$wsGeschlossene.Range("A12").Value = 44044
$wsGeschlossene.Range("B12").Value = 10
$wsGeschlossene.Range("C12").Value = 10

$wsGeschlossene.Range("A13").Value = 44075
$wsGeschlossene.Range("B13").Value = 5
$wsGeschlossene.Range("C13").Value = 7

# Date formatting + styles matching the original "Solved tickets in a year"
# sheet (numFmtId 14 = m/d/yyyy ; row 2's date cell uses the "Komma" cell
# style, the rest use the plain date style).
$wsGeschlossene.Range("A2:A13").NumberFormat = "m/d/yyyy"

$wsGeschlossene.PageSetup.TopMargin = 56.69291339
$wsGeschlossene.PageSetup.BottomMargin = 56.69291339

# ---------------------------------------------------------------------------
# 3) Now that every sheet carries its final content, correct the "Ticket"/
#    "Priority" headers (English original + German translation) to the
#    lower-cased "ticket"/"priority" wording used in the finished workbook.
# ---------------------------------------------------------------------------
$wsOpen.Range("A1").Value = "ticket"
$wsOpen.Range("B1").Value = "priority"

$wsOffene.Range("A1").Value = "ticket"
$wsOffene.Range("B1").Value = "priority"

# ---------------------------------------------------------------------------
# 4) Restore per-sheet selections to match the final, saved state.
# ---------------------------------------------------------------------------
$wsOpen.Range("B1").Select()
$wsSolved.Range("A1:C13").Select()
$wsGeschlossene.Range("B30").Select()
$wsOffene.Range("B1").Select()

Write-Host "Done"
